$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.65176033333333
$ws.Range("H2").Value = 61.955281
$ws.Range("I2").Value = 0.8502439951095444
$ws.Range("J2").Value = 0.8502439951095444
$ws.Range("M2").Value = 58.95713633333333
$ws.Range("N2").Value = 176.871409
$ws.Range("O2").Value = 0.4863146960083892
$ws.Range("P2").Value = 0.4863146960083893
$ws.Range("Q2").Value = 1217.568649495659
$ws.Range("R2").Value = 10958.11784546093
$ws.Range("S2").Value = 0.4134861500146564
$ws.Range("T2").Value = 0.4134861500146565
$ws.Range("G3").Value = 20.65176033333333
$ws.Range("H3").Value = 61.955281
$ws.Range("I3").Value = 0.8502439951095444
$ws.Range("J3").Value = 0.8502439951095444
$ws.Range("O3").Value = 0.07416766570679004
$ws.Range("P3").Value = 0.07416766570679005
$ws.Range("Q3").Value = 185.6909225899766
$ws.Range("R3").Value = 1671.21830330979
$ws.Range("S3").Value = 0.06306061239849031
$ws.Range("T3").Value = 0.06306061239849033
$ws.Range("G4").Value = 20.65176033333333
$ws.Range("H4").Value = 61.955281
$ws.Range("I4").Value = 0.8502439951095444
$ws.Range("J4").Value = 0.8502439951095444
$ws.Range("M4").Value = 42.51661933333333
$ws.Range("N4").Value = 127.549858
$ws.Range("O4").Value = 0.3507032073181665
$ws.Range("P4").Value = 0.3507032073181665
$ws.Range("Q4").Value = 878.0430326555663
$ws.Range("R4").Value = 7902.387293900098
$ws.Range("S4").Value = 0.2981832960879287
$ws.Range("T4").Value = 0.2981832960879287
$ws.Range("G5").Value = 20.65176033333333
$ws.Range("H5").Value = 61.955281
$ws.Range("I5").Value = 0.8502439951095444
$ws.Range("J5").Value = 0.8502439951095444
$ws.Range("M5").Value = 10.76719366666667
$ws.Range("N5").Value = 32.301581
$ws.Range("O5").Value = 0.0888144309666542
$ws.Range("P5").Value = 0.08881443096665421
$ws.Range("Q5").Value = 222.3615030665845
$ws.Range("R5").Value = 2001.253527599261
$ws.Range("S5").Value = 0.0755139366084689
$ws.Range("T5").Value = 0.07551393660846892
$ws.Range("G6").Value = 0.3109896666666667
$ws.Range("H6").Value = 0.932969
$ws.Range("I6").Value = 0.01280361055699766
$ws.Range("J6").Value = 0.01280361055699766
$ws.Range("M6").Value = 58.95713633333333
$ws.Range("N6").Value = 176.871409
$ws.Range("O6").Value = 0.4863146960083892
$ws.Range("P6").Value = 0.4863146960083893
$ws.Range("Q6").Value = 18.33506017592455
$ws.Range("R6").Value = 165.015541583321
$ws.Range("S6").Value = 0.006226583975836121
$ws.Range("T6").Value = 0.006226583975836122
$ws.Range("G7").Value = 0.3109896666666667
$ws.Range("H7").Value = 0.932969
$ws.Range("I7").Value = 0.01280361055699766
$ws.Range("J7").Value = 0.01280361055699766
$ws.Range("O7").Value = 0.07416766570679004
$ws.Range("P7").Value = 0.07416766570679005
$ws.Range("Q7").Value = 2.796272917523333
$ws.Range("R7").Value = 25.16645625771
$ws.Range("S7").Value = 0.0009496139076313303
$ws.Range("T7").Value = 0.0009496139076313305
$ws.Range("G8").Value = 0.3109896666666667
$ws.Range("H8").Value = 0.932969
$ws.Range("I8").Value = 0.01280361055699766
$ws.Range("J8").Value = 0.01280361055699766
$ws.Range("M8").Value = 42.51661933333333
$ws.Range("N8").Value = 127.549858
$ws.Range("O8").Value = 0.3507032073181665
$ws.Range("P8").Value = 0.3507032073181665
$ws.Range("Q8").Value = 13.22222927426689
$ws.Range("R8").Value = 119.000063468402
$ws.Range("S8").Value = 0.004490267287591816
$ws.Range("T8").Value = 0.004490267287591816
$ws.Range("G9").Value = 0.3109896666666667
$ws.Range("H9").Value = 0.932969
$ws.Range("I9").Value = 0.01280361055699766
$ws.Range("J9").Value = 0.01280361055699766
$ws.Range("M9").Value = 10.76719366666667
$ws.Range("N9").Value = 32.301581
$ws.Range("O9").Value = 0.0888144309666542
$ws.Range("P9").Value = 0.08881443096665421
$ws.Range("Q9").Value = 3.348485969332111
$ws.Range("R9").Value = 30.136373723989
$ws.Range("S9").Value = 0.001137145385938394
$ws.Range("T9").Value = 0.001137145385938394
$ws.Range("G10").Value = 3.326466333333334
$ws.Range("H10").Value = 9.979399000000001
$ws.Range("I10").Value = 0.1369523943334579
$ws.Range("J10").Value = 0.1369523943334579
$ws.Range("M10").Value = 58.95713633333333
$ws.Range("N10").Value = 176.871409
$ws.Range("O10").Value = 0.4863146960083892
$ws.Range("P10").Value = 0.4863146960083893
$ws.Range("Q10").Value = 196.1189291225768
$ws.Range("R10").Value = 1765.070362103191
$ws.Range("S10").Value = 0.06660196201789664
$ws.Range("T10").Value = 0.06660196201789664
$ws.Range("G11").Value = 3.326466333333334
$ws.Range("H11").Value = 9.979399000000001
$ws.Range("I11").Value = 0.1369523943334579
$ws.Range("J11").Value = 0.1369523943334579
$ws.Range("O11").Value = 0.07416766570679004
$ws.Range("P11").Value = 0.07416766570679005
$ws.Range("Q11").Value = 29.91002183015667
$ws.Range("R11").Value = 269.19019647141
$ws.Range("S11").Value = 0.01015743940066839
$ws.Range("T11").Value = 0.0101574394006684
$ws.Range("G12").Value = 3.326466333333334
$ws.Range("H12").Value = 9.979399000000001
$ws.Range("I12").Value = 0.1369523943334579
$ws.Range("J12").Value = 0.1369523943334579
$ws.Range("M12").Value = 42.51661933333333
$ws.Range("N12").Value = 127.549858
$ws.Range("O12").Value = 0.3507032073181665
$ws.Range("P12").Value = 0.3507032073181665
$ws.Range("Q12").Value = 141.4301028194824
$ws.Range("R12").Value = 1272.870925375342
$ws.Range("S12").Value = 0.04802964394264598
$ws.Range("T12").Value = 0.04802964394264599
$ws.Range("G13").Value = 3.326466333333334
$ws.Range("H13").Value = 9.979399000000001
$ws.Range("I13").Value = 0.1369523943334579
$ws.Range("J13").Value = 0.1369523943334579
$ws.Range("M13").Value = 10.76719366666667
$ws.Range("N13").Value = 32.301581
$ws.Range("O13").Value = 0.0888144309666542
$ws.Range("P13").Value = 0.08881443096665421
$ws.Range("Q13").Value = 35.81670723664656
$ws.Range("R13").Value = 322.350365129819
$ws.Range("S13").Value = 0.0121633489722469
$ws.Range("T13").Value = 0.01216334897224691